$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.871.20'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  +2.84%  '
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.560.66'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  +2.01%  '
$style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.08%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '615.35'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +6.83%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.04'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('E7').Value = '  +2.69%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.553.89'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('E9').Value = '  +0.00%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.197'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  +5.62%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.29'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +12.89%  '
$ws.Range('E12').Value = '  +1.78%  '
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.71'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('E14').Value = '  +2.54%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.133.80'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +2.01%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.41'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -0.63%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '620.76'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +0.79%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.975.29'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  +3.09%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.563.38'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +2.45%  '
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('E22').Value = '  +0.58%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.48'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -14.07%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.74'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +0.01%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.95'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  +0.45%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.56'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +3.57%  '
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('E31').Value = '  +1.17%  '
$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.07'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('E33').Value = '  -0.03%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.98'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +1.10%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '575.07'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -8.37%  '
$ws.Range('E36').Value = '  +6.50%  '
$ws.Range('E37').Value = '  -0.61%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.87'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +1.97%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '57.62'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0475'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +7.51%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  +6.00%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.374.81'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +1.17%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.320'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -1.13%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.02'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +9.61%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '33.02'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +1.72%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₃0705'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +2.81%  '
$ws.Range('E48').Value = '  +3.20%  '
$ws.Range('E49').Value = '  +1.60%  '
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.79'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +1.66%  '
$ws.Range('E51').Value = '  +1.61%  '
